# Update column G (the "K" column, formerly "Strike#") with newly
# regenerated values for rows 2-12.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @(2, 2, 1, 4, 3, 1, 3, 1, 3, 3, 1)

for ($i = 0; $i -lt $values.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 7).Value = $values[$i]
}
